# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the regenerated data as of commit 456a3b4.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 6509
    $ws.Range("F3").Value = 36
    $ws.Range("F4").Value = 189
    $ws.Range("F6").Value = 121
}
